$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# 1) Update the "总计" (summary) sheet: insert a new top row for 2022-Q3
#    and renumber the existing rows' running index (column A).
# ----------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

$wsTotal.Rows.Item(2).Insert()

$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 2).Value = "2022-Q3"
$wsTotal.Cells.Item(2, 3).Value = 2
$wsTotal.Cells.Item(2, 4).Value = 1.82

# Match formatting: column A carries the bordered/bold style (same as the
# rows below), columns B:D stay unstyled like the rest of the data rows.
$wsTotal.Cells.Item(3, 1).Copy()
$wsTotal.Cells.Item(2, 1).PasteSpecial(-4122)
$wsTotal.Range("B2:D2").ClearFormats()

# Renumber the running index in column A for the rows that shifted down.
$wsTotal.Cells.Item(3, 1).Value = 1
$wsTotal.Cells.Item(4, 1).Value = 2
$wsTotal.Cells.Item(5, 1).Value = 3
$wsTotal.Cells.Item(6, 1).Value = 4

# ----------------------------------------------------------------------
# 2) Add the new "2022-Q3" sheet (placed right after "总计", before
#    "2022-Q1"), cloning the layout/styling of the existing "2022-Q1"
#    sheet and filling in the new quarter's figures.
# ----------------------------------------------------------------------
$wsQ1 = $wb.Worksheets.Item("2022-Q1")
$wsQ1.Copy($wsQ1)
$wsQ3 = $wb.Worksheets.Item("2022-Q1 (2)")
$wsQ3.Name = "2022-Q3"

$wsQ3.Cells.Item(2, 2).Value = "'000906"
$wsQ3.Cells.Item(2, 3).Value = "广发全球精选股票（QDII）美元现汇"
$wsQ3.Cells.Item(2, 4).Value = "'21.88"
$wsQ3.Cells.Item(2, 5).Value = "'79.27"
$wsQ3.Cells.Item(2, 6).Value = "'4.17"
$wsQ3.Cells.Item(2, 7).Value = "'0.9124"
$wsQ3.Cells.Item(2, 8).Value = 7

$wsQ3.Cells.Item(3, 2).Value = "'270023"
$wsQ3.Cells.Item(3, 3).Value = "广发全球精选股票（QDII）"
$wsQ3.Cells.Item(3, 4).Value = "'21.88"
$wsQ3.Cells.Item(3, 5).Value = "'79.27"
$wsQ3.Cells.Item(3, 6).Value = "'4.17"
$wsQ3.Cells.Item(3, 7).Value = "'0.9124"
$wsQ3.Cells.Item(3, 8).Value = 7
